$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped
# from 45175 (2023-09-06) to 45177 (2023-09-08) for every data row
# (rows 2 through 331).
$ws.Range("C2:C331").Value = 45177
